$wb = $excel.ActiveWorkbook

# New row (row 38) data per worksheet, in column order A..I
$rows = @(
    @{
        Sheet = 1
        A = "2025-03-05 21:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    },
    @{
        Sheet = 2
        A = "2025-03-05 21:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    },
    @{
        Sheet = 3
        A = "2025-03-05 21:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    },
    @{
        Sheet = 4
        A = "2025-03-05 21:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
)

foreach ($row in $rows) {
    $ws = $wb.Worksheets.Item($row.Sheet)
    $targetRow = 38

    $ws.Cells.Item($targetRow, 1).Value = $row.A
    $ws.Cells.Item($targetRow, 2).Value = $row.B
    $ws.Cells.Item($targetRow, 3).Value = $row.C
    $ws.Cells.Item($targetRow, 4).Value = $row.D
    $ws.Cells.Item($targetRow, 5).Value = $row.E
    $ws.Cells.Item($targetRow, 6).Value = $row.F

    # Column G holds a 24-digit numeric-looking string that must stay text
    # (mirrors the existing ID_DEC cells) - force text format, assign, then
    # clear the formatting stamp so no stray style index is left behind.
    $gCell = $ws.Cells.Item($targetRow, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $row.G
    $gCell.ClearFormats()

    $ws.Cells.Item($targetRow, 8).Value = $row.H
    $ws.Cells.Item($targetRow, 9).Value = $row.I
}

Write-Output "Appended row 38 to sheets 1-4"
